# IAM.xlsx - "Test Cases" sheet: Runmode column (D) for rows 3-26 changes
# from "Y" to "N" (row 2 stays "Y"). Also update the sheet's active
# selection/scroll position to reflect the new view (B4 top-left,
# D3:D26 selected with D3 as the active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Flip Runmode from Y to N for rows 3 through 26 (row 2 is left as "Y").
$ws.Range("D3:D26").Value = "N"

# Update the view: scroll so B4 is the top-left visible cell, and select
# D3:D26 with D3 as the active cell.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D3:D26").Select()
